$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 959.5
$ws.Range("I15").Value = 959.5
$ws.Range("K15").Value = 2878.5
$ws.Range("M15").Value = -2709.5
$ws.Range("H17").Value = 765.873
$ws.Range("I17").Value = 1305
$ws.Range("J17").Value = 664.15094
$ws.Range("K17").Value = 3915
$ws.Range("L17").Value = 1992.45282
$ws.Range("M17").Value = -3747
$ws.Range("N17").Value = -2328.45282
$ws.Range("H32").Value = 486.8889
$ws.Range("I32").Value = 525
$ws.Range("J32").Value = 456.4
$ws.Range("K32").Value = 525
$ws.Range("L32").Value = 456.4
$ws.Range("M32").Value = -199
$ws.Range("N32").Value = -1108.4
$ws.Range("H123").Value = 41835
$ws.Range("J123").Value = 41835
$ws.Range("L123").Value = 41835
$ws.Range("N123").Value = -51635
$ws.Range("H130").Value = 41835
$ws.Range("J130").Value = 41835
$ws.Range("L130").Value = 41835
$ws.Range("N130").Value = -51875
$ws.Range("H137").Value = 2370.6829
$ws.Range("I137").Value = 1163.2812
$ws.Range("J137").Value = 6663.6665
$ws.Range("K137").Value = 3489.8436
$ws.Range("L137").Value = 19990.9995
$ws.Range("M137").Value = -939.8435999999997
$ws.Range("N137").Value = -25090.9995
$ws.Range("H138").Value = 3866.2
$ws.Range("I138").Value = 829
$ws.Range("K138").Value = 2487
$ws.Range("M138").Value = 2653

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1069.1428
$ws.Range("I61").Value = 883.64703
$ws.Range("J61").Value = 1355.8182
$ws.Range("K61").Value = 883.64703
$ws.Range("L61").Value = 1355.8182
$ws.Range("M61").Value = -671.64703
$ws.Range("N61").Value = -1779.8182
$ws.Range("H74").Value = 2544.4285
$ws.Range("I74").Value = 2374.9524
$ws.Range("J74").Value = 3561.2856
$ws.Range("K74").Value = 2374.9524
$ws.Range("L74").Value = 3561.2856
$ws.Range("M74").Value = -1500.9524
$ws.Range("N74").Value = -5309.2856
$ws.Range("H77").Value = 2544.4285
$ws.Range("I77").Value = 2374.9524
$ws.Range("J77").Value = 3561.2856
$ws.Range("K77").Value = 11874.762
$ws.Range("L77").Value = 17806.428
$ws.Range("M77").Value = -7506.762000000001
$ws.Range("N77").Value = -26542.428
$ws.Range("H88").Value = 22225554
$ws.Range("I88").Value = 33335832
$ws.Range("K88").Value = 33335832
$ws.Range("M88").Value = -33335426
$ws.Range("H91").Value = 22225554
$ws.Range("I91").Value = 33335832
$ws.Range("K91").Value = 33335832
$ws.Range("M91").Value = -33334428
$ws.Range("H132").Value = 1882.3922
$ws.Range("I132").Value = 1229.7
$ws.Range("J132").Value = 4255.8184
$ws.Range("K132").Value = 3689.1
$ws.Range("L132").Value = 12767.4552
$ws.Range("M132").Value = -1159.1
$ws.Range("N132").Value = -17827.4552
$ws.Range("H136").Value = 1069.1428
$ws.Range("I136").Value = 883.64703
$ws.Range("J136").Value = 1355.8182
$ws.Range("K136").Value = 2650.94109
$ws.Range("L136").Value = 4067.4546
$ws.Range("M136").Value = -100.9410899999998
$ws.Range("N136").Value = -9167.454600000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2458.1428
$ws.Range("I86").Value = 2250
$ws.Range("J86").Value = 2735.6667
$ws.Range("K86").Value = 2250
$ws.Range("L86").Value = 2735.6667
$ws.Range("M86").Value = -1127
$ws.Range("N86").Value = -4981.6667
$ws.Range("H89").Value = 2458.1428
$ws.Range("I89").Value = 2250
$ws.Range("J89").Value = 2735.6667
$ws.Range("K89").Value = 11250
$ws.Range("L89").Value = 13678.3335
$ws.Range("M89").Value = -5634
$ws.Range("N89").Value = -24910.3335
$ws.Range("H134").Value = 2458.7441
$ws.Range("I134").Value = 1244.625
$ws.Range("K134").Value = 3733.875
$ws.Range("M134").Value = -1198.875

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 12823133
$ws.Range("I31").Value = 1431.963
$ws.Range("J31").Value = 41671960
$ws.Range("K31").Value = 1431.963
$ws.Range("L31").Value = 41671960
$ws.Range("M31").Value = -1136.963
$ws.Range("N31").Value = -41672550
$ws.Range("H34").Value = 12823133
$ws.Range("I34").Value = 1431.963
$ws.Range("J34").Value = 41671960
$ws.Range("K34").Value = 1431.963
$ws.Range("L34").Value = 41671960
$ws.Range("M34").Value = -1229.963
$ws.Range("N34").Value = -41672364
$ws.Range("H58").Value = 2042.6461
$ws.Range("I58").Value = 1766.2931
$ws.Range("J58").Value = 4332.4287
$ws.Range("K58").Value = 1766.2931
$ws.Range("L58").Value = 4332.4287
$ws.Range("M58").Value = -1563.2931
$ws.Range("N58").Value = -4738.4287
$ws.Range("H99").Value = 14292085
$ws.Range("I99").Value = 22226722
$ws.Range("J99").Value = 9738
$ws.Range("K99").Value = 22226722
$ws.Range("L99").Value = 9738
$ws.Range("M99").Value = -22225224
$ws.Range("N99").Value = -12734
$ws.Range("H126").Value = 14292085
$ws.Range("I126").Value = 22226722
$ws.Range("J126").Value = 9738
$ws.Range("K126").Value = 66680166
$ws.Range("L126").Value = 29214
$ws.Range("M126").Value = -66677696
$ws.Range("N126").Value = -34154
$ws.Range("H132").Value = 3108.2222
$ws.Range("I132").Value = 2738.6785
$ws.Range("J132").Value = 4401.625
$ws.Range("K132").Value = 8216.0355
$ws.Range("L132").Value = 13204.875
$ws.Range("M132").Value = -5686.0355
$ws.Range("N132").Value = -18264.875
$ws.Range("H136").Value = 2042.6461
$ws.Range("I136").Value = 1766.2931
$ws.Range("J136").Value = 4332.4287
$ws.Range("K136").Value = 5298.879300000001
$ws.Range("L136").Value = 12997.2861
$ws.Range("M136").Value = -2748.879300000001
$ws.Range("N136").Value = -18097.2861

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H137").Value = 2361.907
$ws.Range("I137").Value = 662.5333000000001
$ws.Range("J137").Value = 3272.2856
$ws.Range("K137").Value = 1987.5999
$ws.Range("L137").Value = 9816.856800000001
$ws.Range("M137").Value = 3112.4001
$ws.Range("N137").Value = -20016.8568

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H123").Value = 11985.77
$ws.Range("J123").Value = 11985.77
$ws.Range("L123").Value = 11985.77
$ws.Range("N123").Value = -16885.77

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 26321852
$ws.Range("I40").Value = 50004372
$ws.Range("K40").Value = 50004372
$ws.Range("M40").Value = -50004236
$ws.Range("H122").Value = 4624.55
$ws.Range("I122").Value = 2623.9167
$ws.Range("K122").Value = 7871.750100000001
$ws.Range("M122").Value = -5421.750100000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 436.1905
$ws.Range("I113").Value = 326.5
$ws.Range("J113").Value = 655.5714
$ws.Range("K113").Value = 979.5
$ws.Range("L113").Value = 1966.7142
$ws.Range("M113").Value = 1190.5
$ws.Range("N113").Value = -6306.7142
$ws.Range("H122").Value = 6012.1113
$ws.Range("I122").Value = 3184.8333
$ws.Range("K122").Value = 9554.499899999999
$ws.Range("M122").Value = -7104.499899999999
$ws.Range("H132").Value = 7577658
$ws.Range("I132").Value = 1537.3334
$ws.Range("J132").Value = 19610320
$ws.Range("K132").Value = 4612.0002
$ws.Range("L132").Value = 58830960
$ws.Range("M132").Value = -2082.0002
$ws.Range("N132").Value = -58836020
